# mudando layout do dash
# Add 6 new transaction rows (87-92) to the bottom of the sheet, following
# the same layout/pattern as the existing rows (A=data/date, B=hora/time
# text, C=preco/price, D=site, E=cor/color).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 87; Date = 45228; Time = "10:00"; Price = 3177.32; Site = "amazon";         Color = "preto" },
    @{ Row = 88; Date = 45228; Time = "10:00"; Price = 2635;    Site = "mercado livre";  Color = "preto" },
    @{ Row = 89; Date = 45229; Time = "20:58"; Price = 3176.33; Site = "amazon";         Color = "preto" },
    @{ Row = 90; Date = 45229; Time = "20:58"; Price = 2599;    Site = "mercado livre";  Color = "preto" },
    @{ Row = 91; Date = 45230; Time = "12:40"; Price = 3176.33; Site = "amazon";         Color = "preto" },
    @{ Row = 92; Date = 45230; Time = "12:40"; Price = 2599;    Site = "mercado livre";  Color = "preto" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $cellA = $ws.Cells.Item($row, 1)
    $cellA.Value2 = $r.Date
    $cellA.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value2 = $r.Time
    $ws.Cells.Item($row, 3).Value2 = $r.Price
    $ws.Cells.Item($row, 4).Value2 = $r.Site
    $ws.Cells.Item($row, 5).Value2 = $r.Color
}
